# Atualizacao de bases das ligas - apply cell-level corrections to
# match, team, and odds data as described by the commit diff.
# The changes affect two groups of rows:
#  - Rows 134-145 (except 138, 141, which are unchanged): the home/away
#    team assignment (and every odds-related column tied to it) was
#    corrected, which moves a differing set of values into each row.
#  - Rows 200, 202, 203, 204: pure odds-value corrections (no team swap).
#
# Each entry below is (row, column index, new value) matching the exact
# target state described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(134, 2, 7483306),
    @(134, 6, "Tecnico Universitario"),
    @(134, 7, "Club Atletico Libertad"),
    @(134, 8, 1),
    @(134, 10, "D"),
    @(134, 11, 1.5),
    @(134, 12, 4.333),
    @(134, 13, 5.75),
    @(134, 14, 1.533),
    @(134, 15, 4.2),
    @(134, 16, 5.5),
    @(134, 17, -1),
    @(134, 18, 1.925),
    @(134, 19, 1.875),
    @(134, 20, 2.25),
    @(134, 21, 1.8),
    @(134, 22, 2),
    @(134, 23, -1),
    @(134, 24, 3.2),
    @(134, 26, -1),
    @(134, 27, 0.875),
    @(134, 28, -0.5),
    @(134, 29, 0.5),
    @(135, 2, 7482832),
    @(135, 6, "Barcelona Guayaquil"),
    @(135, 7, "Guayaquil City"),
    @(135, 8, 2),
    @(135, 10, "H"),
    @(135, 11, 1.363),
    @(135, 12, 5),
    @(135, 13, 7.5),
    @(135, 14, 1.444),
    @(135, 15, 4),
    @(135, 16, 8),
    @(135, 17, -1.25),
    @(135, 18, 2.05),
    @(135, 19, 1.75),
    @(135, 20, 2.5),
    @(135, 21, 1.95),
    @(135, 22, 1.85),
    @(135, 23, 0.444),
    @(135, 24, -1),
    @(135, 26, -0.5),
    @(135, 27, 0.375),
    @(135, 28, 0.95),
    @(135, 29, -1),
    @(136, 2, 7483188),
    @(136, 6, "Gualaceo SC"),
    @(136, 7, "Emelec"),
    @(136, 8, 0),
    @(136, 11, 3.6),
    @(136, 12, 3.3),
    @(136, 13, 2.05),
    @(136, 14, 2.6),
    @(136, 15, 3.25),
    @(136, 16, 2.75),
    @(136, 17, 0),
    @(136, 18, 1.8),
    @(136, 19, 2),
    @(136, 21, 1.975),
    @(136, 22, 1.825),
    @(136, 25, 1.75),
    @(136, 26, -1),
    @(136, 27, 1),
    @(136, 28, -1),
    @(136, 29, 0.825),
    @(137, 2, 7482867),
    @(137, 6, "Cumbaya FC"),
    @(137, 7, "LDU Quito"),
    @(137, 8, 1),
    @(137, 11, 5.25),
    @(137, 12, 3.75),
    @(137, 13, 1.65),
    @(137, 14, 9),
    @(137, 15, 4.5),
    @(137, 16, 1.363),
    @(137, 17, 1.25),
    @(137, 18, 1.975),
    @(137, 19, 1.825),
    @(137, 21, 1.825),
    @(137, 22, 1.975),
    @(137, 25, 0.363),
    @(137, 26, 0.4875),
    @(137, 27, -0.5),
    @(137, 28, 0.825),
    @(137, 29, -1),
    @(139, 2, 7528859),
    @(139, 6, "Club Atletico Libertad"),
    @(139, 7, "Cumbaya FC"),
    @(139, 8, 3),
    @(139, 9, 1),
    @(139, 10, "H"),
    @(139, 11, 1.727),
    @(139, 13, 4.333),
    @(139, 14, 1.4),
    @(139, 15, 4.2),
    @(139, 16, 7),
    @(139, 17, -1.25),
    @(139, 18, 2),
    @(139, 19, 1.8),
    @(139, 21, 1.95),
    @(139, 22, 1.85),
    @(139, 23, 0.3999999999999999),
    @(139, 25, -1),
    @(139, 26, 1),
    @(139, 27, -1),
    @(139, 28, 0.95),
    @(139, 29, -1),
    @(140, 2, 7528849),
    @(140, 6, "Guayaquil City"),
    @(140, 7, "Gualaceo SC"),
    @(140, 8, 0),
    @(140, 9, 2),
    @(140, 10, "A"),
    @(140, 11, 1.833),
    @(140, 13, 3.75),
    @(140, 14, 2.15),
    @(140, 15, 3.4),
    @(140, 16, 3),
    @(140, 17, -0.25),
    @(140, 18, 1.825),
    @(140, 19, 1.975),
    @(140, 21, 1.85),
    @(140, 22, 1.95),
    @(140, 23, -1),
    @(140, 25, 2),
    @(140, 26, -1),
    @(140, 27, 0.9750000000000001),
    @(140, 28, -1),
    @(140, 29, 0.95),
    @(142, 2, 7528858),
    @(142, 6, "Orense"),
    @(142, 7, "SD Aucas"),
    @(142, 8, 1),
    @(142, 10, "A"),
    @(142, 11, 2.2),
    @(142, 12, 3.2),
    @(142, 13, 3.2),
    @(142, 14, 1.95),
    @(142, 15, 3.2),
    @(142, 16, 3.8),
    @(142, 17, -0.5),
    @(142, 18, 1.95),
    @(142, 19, 1.85),
    @(142, 21, 1.85),
    @(142, 22, 1.95),
    @(142, 24, -1),
    @(142, 25, 2.8),
    @(142, 26, -1),
    @(142, 27, 0.8500000000000001),
    @(142, 28, 0.8500000000000001),
    @(143, 2, 7528857),
    @(143, 6, "Universidad Catolica del Ecuador"),
    @(143, 7, "Barcelona Guayaquil"),
    @(143, 8, 0),
    @(143, 10, "A"),
    @(143, 11, 1.533),
    @(143, 12, 4),
    @(143, 13, 5.5),
    @(143, 14, 1.5),
    @(143, 15, 4.333),
    @(143, 16, 5.25),
    @(143, 17, -1),
    @(143, 18, 1.8),
    @(143, 19, 2),
    @(143, 20, 3),
    @(143, 21, 1.975),
    @(143, 22, 1.825),
    @(143, 23, -1),
    @(143, 25, 4.25),
    @(143, 26, -1),
    @(143, 27, 1),
    @(143, 28, -1),
    @(143, 29, 0.825),
    @(144, 2, 7528848),
    @(144, 6, "Emelec"),
    @(144, 7, "Deportivo Cuenca"),
    @(144, 8, 2),
    @(144, 10, "H"),
    @(144, 11, 1.75),
    @(144, 12, 3.5),
    @(144, 13, 4.2),
    @(144, 14, 2.4),
    @(144, 15, 3.1),
    @(144, 16, 2.75),
    @(144, 17, -0.25),
    @(144, 18, 2.05),
    @(144, 19, 1.75),
    @(144, 20, 2.25),
    @(144, 21, 1.8),
    @(144, 22, 2),
    @(144, 23, 1.4),
    @(144, 25, -1),
    @(144, 26, 1.05),
    @(144, 27, -1),
    @(144, 28, 0.8),
    @(144, 29, -1),
    @(145, 2, 7528852),
    @(145, 6, "Delfin SC"),
    @(145, 7, "Tecnico Universitario"),
    @(145, 8, 2),
    @(145, 10, "D"),
    @(145, 11, 2.1),
    @(145, 12, 3.4),
    @(145, 13, 3.1),
    @(145, 14, 2.1),
    @(145, 15, 3.4),
    @(145, 16, 3.1),
    @(145, 17, -0.25),
    @(145, 18, 1.8),
    @(145, 19, 2),
    @(145, 21, 1.9),
    @(145, 22, 1.9),
    @(145, 24, 2.4),
    @(145, 25, -1),
    @(145, 26, -0.5),
    @(145, 27, 0.5),
    @(145, 28, 0.8999999999999999),
    @(200, 14, 3.1),
    @(200, 15, 3.1),
    @(200, 16, 2.2),
    @(200, 18, 1.85),
    @(200, 19, 1.95),
    @(200, 21, 1.95),
    @(200, 22, 1.85),
    @(202, 18, 1.825),
    @(202, 19, 1.975),
    @(203, 14, 1.222),
    @(203, 15, 5.75),
    @(203, 16, 11),
    @(203, 17, -1.75),
    @(203, 18, 1.95),
    @(203, 19, 1.85),
    @(203, 21, 1.9),
    @(203, 22, 1.9),
    @(204, 14, 2.6),
    @(204, 16, 2.6),
    @(204, 18, 1.875),
    @(204, 19, 1.925)
)

foreach ($change in $changes) {
    $row = $change[0]
    $col = $change[1]
    $val = $change[2]
    $ws.Cells.Item($row, $col).Value = $val
}

Write-Output "Applied $($changes.Count) cell updates"
